# Adds two new "Title and Content" slides (indexes 18 and 19) to the end of
# the deck: "Team Members - Physicians" and "Team Members - Support Staff".
# Layout 2 on the master corresponds to ppt/slideLayouts/slideLayout2.xml
# ("Title and Content"), the same layout used by most of the existing
# content slides (e.g. slide3.xml).

$p = $ppt.ActivePresentation

function Add-BulletlessParagraph($textRange, [int]$index) {
    $para = $textRange.Paragraphs($index, 1)
    $para.ParagraphFormat.Bullet.Visible = 0
    $para.IndentLevel = 0
}

function Add-BulletedParagraph($textRange, [int]$index) {
    $para = $textRange.Paragraphs($index, 1)
    $para.IndentLevel = 0
}

# ---------------------------------------------------------------------
# Slide 18 - Team Members - Physicians
# ---------------------------------------------------------------------
$s18 = $p.Slides.Add(18, 2)

$title18 = $s18.Shapes.Item(1).TextFrame.TextRange
$title18.Text = "Team Members - Physicians"
Add-BulletlessParagraph $title18 1

$body18 = $s18.Shapes.Item(2).TextFrame.TextRange
$lines18 = @(
    "Primary Care Provider",
    "Gastroenterologist",
    "Medical Oncologist (chemotherapy)",
    "Radiation Oncologist (radiation)",
    "Surgeons",
    "Jonathan Salo MD",
    "Jeffrey Hagen MD",
    "Michael Roach MD"
)
$body18.Text = [string]::Join("`r", $lines18)

for ($i = 1; $i -le 5; $i++) {
    Add-BulletlessParagraph $body18 $i
}
for ($i = 6; $i -le 8; $i++) {
    Add-BulletedParagraph $body18 $i
}

# ---------------------------------------------------------------------
# Slide 19 - Team Members - Support Staff
# ---------------------------------------------------------------------
$s19 = $p.Slides.Add(19, 2)

$title19 = $s19.Shapes.Item(1).TextFrame.TextRange
$title19.Text = "Team Members - Support Staff"
Add-BulletlessParagraph $title19 1

$body19 = $s19.Shapes.Item(2).TextFrame.TextRange
$lines19 = @(
    "Dietitian - Liz Koch",
    "Nurses",
    "Matthew Carpenter RN",
    "Brandon Galloway LPN",
    "Navigator - Laura Swift"
)
$body19.Text = [string]::Join("`r", $lines19)

Add-BulletlessParagraph $body19 1
Add-BulletlessParagraph $body19 2
Add-BulletedParagraph   $body19 3
Add-BulletedParagraph   $body19 4
Add-BulletlessParagraph $body19 5
